$wb = $excel.ActiveWorkbook

# --- Update status text "Ready for handoff" -> "In Translation" ---
# This shared string is referenced from:
#   Overview sheet: E2, F2
#   zh-cn sheet:     C2
#   de-de sheet:     C2
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "In Translation"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "In Translation"

# --- Narrow the "zh-cn"/"de-de" status columns ---
# Overview!E:F and the Status column (C) on each language sheet shrink from
# ~17.22 chars to ~13.41 chars. Excel's ColumnWidth setter snaps to whole
# pixel increments, so we choose the input value whose rounded result lands
# on the pixel width closest to the target.
$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5

$wsZhCn.Columns.Item(3).ColumnWidth = 12.5
$wsDeDe.Columns.Item(3).ColumnWidth = 12.5
